# Generate Report for Handback
# Update the timestamp values that reflect when the handback report was generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first row
$wsOverview.Range("G2").Value = "2016-09-07 05:25:15"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsZhCn.Range("H2").Value = "2016-09-07 05:25:03"
$wsZhCn.Range("K2").Value = "2016-09-07 05:25:33"

# de-de sheet: Correspond Handback DateTime for first row
$wsDeDe.Range("K2").Value = "2016-09-07 05:25:41"
